$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column cells are treated as literal text (matching source data
# which stores prices/volumes as text strings, not numbers), so values such as
# "0.999" or "4.70" are preserved exactly rather than being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "74.475.78"
$ws.Range("E2").Value = "  +5.89%  "

# Row 3
$ws.Range("D3").Value = "2.650.95"
$ws.Range("E3").Value = "  +7.34%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "185.69"
$ws.Range("E5").Value = "  +9.81%  "

# Row 6
$ws.Range("D6").Value = "583.89"
$ws.Range("E6").Value = "  +2.59%  "

# Row 7
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.12%  "

# Row 8
$ws.Range("E8").Value = "  +3.26%  "

# Row 9
$ws.Range("E9").Value = "  +8.88%  "

# Row 10
$ws.Range("D10").Value = "2.649.92"
$ws.Range("E10").Value = "  +7.32%  "

# Row 11
$ws.Range("D11").Value = "0.163"
$ws.Range("E11").Value = "  +1.11%  "

# Row 12
$ws.Range("E12").Value = "  +4.88%  "

# Row 13
$ws.Range("D13").Value = "4.70"
$ws.Range("E13").Value = "  -0.53%  "

# Row 14
$ws.Range("D14").Value = "3.137.74"
$ws.Range("E14").Value = "  +7.71%  "

# Row 15
$ws.Range("D15").Value = "74.386.48"
$ws.Range("E15").Value = "  +5.95%  "

# Row 16
$ws.Range("E16").Value = "  +0.60%  "

# Row 17
$ws.Range("D17").Value = "26.35"
$ws.Range("E17").Value = "  +8.27%  "

# Row 18
$ws.Range("D18").Value = "2.653.89"
$ws.Range("E18").Value = "  +7.82%  "

# Row 19
$ws.Range("D19").Value = "9.24"
$ws.Range("E19").Value = "  +27.68%  "

# Row 20
$ws.Range("D20").Value = "11.84"
$ws.Range("E20").Value = "  +8.44%  "

# Row 21
$ws.Range("D21").Value = "370.95"
$ws.Range("E21").Value = "  +7.58%  "

# Row 22
$ws.Range("D22").Value = "2.27"
$ws.Range("E22").Value = "  +11.94%  "

# Row 23
$ws.Range("D23").Value = "4.06"
$ws.Range("E23").Value = "  +3.81%  "

# Row 24
$ws.Range("E24").Value = "  +2.80%  "

# Row 25
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.14%  "

# Row 26
$ws.Range("D26").Value = "69.43"
$ws.Range("E26").Value = "  +4.00%  "

# Row 27
$ws.Range("D27").Value = "4.10"
$ws.Range("E27").Value = "  +4.47%  "

# Row 28
$ws.Range("D28").Value = "9.33"
$ws.Range("E28").Value = "  +8.27%  "

# Row 29
$ws.Range("D29").Value = "2.788.85"
$ws.Range("E29").Value = "  +7.85%  "

# Row 30
$ws.Range("E30").Value = "  +15.11%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0933"
$ws.Range("E31").Value = "  +6.92%  "

# Row 32
$ws.Range("D32").Value = "517.53"
$ws.Range("E32").Value = "  +12.83%  "

# Row 33
$ws.Range("D33").Value = "1.40"
$ws.Range("E33").Value = "  +10.84%  "

# Row 34
$ws.Range("E34").Value = "  +2.98%  "

# Row 35
$ws.Range("D35").Value = "1.74"
$ws.Range("E35").Value = "  +6.06%  "

# Row 36
$ws.Range("E36").Value = "  -0.07%  "

# Row 37
$ws.Range("D37").Value = "162.78"
$ws.Range("E37").Value = "  +0.88%  "

# Row 38
$ws.Range("D38").Value = "0.118"
$ws.Range("E38").Value = "  +3.65%  "

# Row 39
$ws.Range("D39").Value = "19.10"
$ws.Range("E39").Value = "  +4.31%  "

# Row 40
$ws.Range("D40").Value = "19.34"
$ws.Range("E40").Value = "  +1.17%  "

# Row 41
$ws.Range("E41").Value = "  +0.02%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "4.96"
$ws.Range("E42").Value = "  +11.02%  "

# Row 43
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "167.94"
$ws.Range("E43").Value = "  +25.15%  "

# Row 44
$ws.Range("E44").Value = "  +6.30%  "

# Row 45
$ws.Range("D45").Value = "1.66"
$ws.Range("E45").Value = "  +6.54%  "

# Row 46
$ws.Range("E46").Value = "  +6.69%  "

# Row 47
$ws.Range("D47").Value = "39.04"
$ws.Range("E47").Value = "  +2.92%  "

# Row 48
$ws.Range("D48").Value = "2.33"
$ws.Range("E48").Value = "  +7.57%  "

# Row 49
$ws.Range("D49").Value = "0.0845"
$ws.Range("E49").Value = "  +16.26%  "

# Row 50
$ws.Range("E50").Value = "  +5.57%  "

# Row 51
$ws.Range("E51").Value = "  +5.86%  "
